$wb = $excel.ActiveWorkbook

# Remove the "Signature" sheet entirely
$wsToDelete = $wb.Worksheets.Item("Signature")
$wsToDelete.Delete()

# Hotzones: rename a few cells (F1, G1, H1)
$wsHotzones = $wb.Worksheets.Item("Hotzones")
$wsHotzones.Range("F1").Value = "Mid Range Left"
$wsHotzones.Range("G1").Value = "Mid-Range Right"
$wsHotzones.Range("H1").Value = "Mid-Range Right-Center"

$wsVitals = $wb.Worksheets.Item("Vitals")
$wsVitals.Range("A1:DV1").ClearContents()
$arrVitals = New-Object 'object[,]' 1,111
$arrVitals[0,0] = "ACTIONSHOTID"
$arrVitals[0,1] = "Age"
$arrVitals[0,2] = "ARM_SCALE"
$arrVitals[0,3] = "AUDIOSIGNATUREID"
$arrVitals[0,4] = "Average Percent"
$arrVitals[0,5] = "Birth Day"
$arrVitals[0,6] = "Birth Month"
$arrVitals[0,7] = "Birth Year"
$arrVitals[0,8] = "BODYLENGTH"
$arrVitals[0,9] = "BODY_SHAPE"
$arrVitals[0,10] = "Boom Percentage"
$arrVitals[0,11] = "BOOM_OR_BUST"
$arrVitals[0,12] = "Bust Percentage"
$arrVitals[0,13] = "CAREERENDINGINJURYALLOWED"
$arrVitals[0,14] = "College"
$arrVitals[0,15] = "CONTRACT_TEAM"
$arrVitals[0,16] = "Current Team"
$arrVitals[0,17] = "Current Team Address"
$arrVitals[0,18] = "Dominant Hand"
$arrVitals[0,19] = "DOMINANTDUNKHAND"
$arrVitals[0,20] = "DRAFTED_TEAM"
$arrVitals[0,21] = "DRAFTPICK"
$arrVitals[0,22] = "DRAFTROUND"
$arrVitals[0,23] = "DRAFTYEAR"
$arrVitals[0,24] = "Dunk Hand"
$arrVitals[0,25] = "EYE_COLOR"
$arrVitals[0,26] = "Face ID"
$arrVitals[0,27] = "Financial Security"
$arrVitals[0,28] = "First Name"
$arrVitals[0,29] = "FIRSTNAMESHOWCASE"
$arrVitals[0,30] = "Force Non Starter"
$arrVitals[0,31] = "G-LEAGUE_TEAM_ID"
$arrVitals[0,32] = "GENDER"
$arrVitals[0,33] = "HADIMPORTANTINJURY"
$arrVitals[0,34] = "HAIR_LENGTH"
$arrVitals[0,35] = "HAND_SCALE"
$arrVitals[0,36] = "Headshot ID"
$arrVitals[0,37] = "Height"
$arrVitals[0,38] = "HOMETOWN_TEAM"
$arrVitals[0,39] = "Injury 1 Duration"
$arrVitals[0,40] = "Injury 1 Type"
$arrVitals[0,41] = "INJURY1BODY"
$arrVitals[0,42] = "INJURY1DAY"
$arrVitals[0,43] = "INJURY1STATUS"
$arrVitals[0,44] = "INJURY2BODY"
$arrVitals[0,45] = "INJURY2DAY"
$arrVitals[0,46] = "INJURY2STATUS"
$arrVitals[0,47] = "INJURY2TYPE"
$arrVitals[0,48] = "IS_ACTIVE"
$arrVitals[0,49] = "IS_DLEAGUE"
$arrVitals[0,50] = "IS_DRAFTED_IN_TEAM_EXPANSION"
$arrVitals[0,51] = "IS_DRAFT_PROSPECT"
$arrVitals[0,52] = "IS_ELIGIBLE_FOR_FANTASY_DRAFT"
$arrVitals[0,53] = "IS_GENERATED"
$arrVitals[0,54] = "IS_GLEAGUE_SEND_DOWN"
$arrVitals[0,55] = "IS_HIDDEN"
$arrVitals[0,56] = "IS_HISTORIC"
$arrVitals[0,57] = "IS_PROTECTED_IN_TEAM_EXTENSION_DRAFT"
$arrVitals[0,58] = "IS_SUMMER_LEAGUE_ATTENDEE"
$arrVitals[0,59] = "Jersey Number"
$arrVitals[0,60] = "Last Name"
$arrVitals[0,61] = "LOWER_SCALE"
$arrVitals[0,62] = "Loyalty"
$arrVitals[0,63] = "MAKE_A_COACH"
$arrVitals[0,64] = "Maximum Potential"
$arrVitals[0,65] = "Minimum Potential"
$arrVitals[0,66] = "MURAL_ID"
$arrVitals[0,67] = "MURAL_TEAM"
$arrVitals[0,68] = "MUST_RETIRE_NEXT_YEAR"
$arrVitals[0,69] = "MYTEAM_DUPLICATE_ID"
$arrVitals[0,70] = "NAMEORDER"
$arrVitals[0,71] = "NAMESID"
$arrVitals[0,72] = "NECK_HEAD_SCALE"
$arrVitals[0,73] = "NICKNAME"
$arrVitals[0,74] = "NICKNAME1"
$arrVitals[0,75] = "NICKNAME_CREATED"
$arrVitals[0,76] = "ORIGINALINJURY1TYPE"
$arrVitals[0,77] = "ORIGINALINJURY2TYPE"
$arrVitals[0,78] = "Peak End"
$arrVitals[0,79] = "Peak Start"
$arrVitals[0,80] = "PERSONALITY"
$arrVitals[0,81] = "PHOTOID"
$arrVitals[0,82] = "Play For Winner"
$arrVitals[0,83] = "Play Initiator"
$arrVitals[0,84] = "Play Type 1"
$arrVitals[0,85] = "Play Type 2"
$arrVitals[0,86] = "Play Type 3"
$arrVitals[0,87] = "Play Type 4"
$arrVitals[0,88] = "Portrait ID"
$arrVitals[0,89] = "Portrait Team 1"
$arrVitals[0,90] = "Portrait Team 2"
$arrVitals[0,91] = "Position"
$arrVitals[0,92] = "PREVIOUS_TEAM"
$arrVitals[0,93] = "QUALIFIESFORWORLDTEAM"
$arrVitals[0,94] = "RECENT_DRAFT_SIGNING"
$arrVitals[0,95] = "RETIREMENT"
$arrVitals[0,96] = "Secondary Position"
$arrVitals[0,97] = "SHOULDERWIDTH"
$arrVitals[0,98] = "SIGNATUREID"
$arrVitals[0,99] = "SKINCOLOR"
$arrVitals[0,100] = "SKINTYPE"
$arrVitals[0,101] = "THIRD_POSITION"
$arrVitals[0,102] = "TYPE"
$arrVitals[0,103] = "UNIQUEID"
$arrVitals[0,104] = "UNIQUESIGNATUREID"
$arrVitals[0,105] = "Used Retirement Grace"
$arrVitals[0,106] = "VOICETYPE"
$arrVitals[0,107] = "WAS_DRAFTED"
$arrVitals[0,108] = "Weight"
$arrVitals[0,109] = "Wingspan"
$arrVitals[0,110] = "Years Pro"
$wsVitals.Range("A1").Resize(1,111).Value = $arrVitals

$wsTend = $wb.Worksheets.Item("Tendencies")
$wsTend.Range("A1:CN1").ClearContents()
$arrTend = New-Object 'object[,]' 1,88
$arrTend[0,0] = "Alley Oop"
$arrTend[0,1] = "Attack Strong On Drive"
$arrTend[0,2] = "Block Shot"
$arrTend[0,3] = "Contest Shot"
$arrTend[0,4] = "Contested Jumper 3pt"
$arrTend[0,5] = "Contested Jumper Mid"
$arrTend[0,6] = "Crash"
$arrTend[0,7] = "Dish To Open Man"
$arrTend[0,8] = "Dribble Spin"
$arrTend[0,9] = "Drive"
$arrTend[0,10] = "Drive Pull Up 3pt"
$arrTend[0,11] = "Drive Pull Up Mid"
$arrTend[0,12] = "Drive Right"
$arrTend[0,13] = "Driving Behind The Back"
$arrTend[0,14] = "Driving Double Crossover"
$arrTend[0,15] = "Driving Dribble Hesitation"
$arrTend[0,16] = "Driving Dunk Tendency"
$arrTend[0,17] = "Driving In And Out"
$arrTend[0,18] = "Driving Layup Tendency"
$arrTend[0,19] = "Driving Step Back"
$arrTend[0,20] = "Euro Step Layup"
$arrTend[0,21] = "Flashy Dunk"
$arrTend[0,22] = "Flashy Pass"
$arrTend[0,23] = "Floater"
$arrTend[0,24] = "Foul"
$arrTend[0,25] = "Hard Foul"
$arrTend[0,26] = "Hop Step Layup"
$arrTend[0,27] = "Iso Vs Average Defender"
$arrTend[0,28] = "Iso Vs Elite Defender"
$arrTend[0,29] = "Iso Vs Good Defender"
$arrTend[0,30] = "Iso Vs Poor Defender"
$arrTend[0,31] = "No Driving Dribble Move"
$arrTend[0,32] = "No Setup Dribble"
$arrTend[0,33] = "Off Screen Drive"
$arrTend[0,34] = "Off Screen Shot 3pt"
$arrTend[0,35] = "Off Screen Shot Mid"
$arrTend[0,36] = "On Ball Steal"
$arrTend[0,37] = "Pass Interception"
$arrTend[0,38] = "Play Discipline"
$arrTend[0,39] = "Post Aggressive Backdown"
$arrTend[0,40] = "Post Back Down"
$arrTend[0,41] = "Post Drive"
$arrTend[0,42] = "Post Drop Step"
$arrTend[0,43] = "Post Fade Left"
$arrTend[0,44] = "Post Fade Right"
$arrTend[0,45] = "Post Hook Left"
$arrTend[0,46] = "Post Hook Right"
$arrTend[0,47] = "Post Hop Shot Tendency"
$arrTend[0,48] = "Post Hop Step"
$arrTend[0,49] = "Post Shimmy Shot"
$arrTend[0,50] = "Post Spin"
$arrTend[0,51] = "Post Step Back Shot"
$arrTend[0,52] = "Post Up"
$arrTend[0,53] = "Post Up And Under"
$arrTend[0,54] = "Putback Dunk"
$arrTend[0,55] = "Roll Vs Pop"
$arrTend[0,56] = "Setup With Hesitation"
$arrTend[0,57] = "Setup With Sizeup"
$arrTend[0,58] = "Shoot"
$arrTend[0,59] = "Shoot From Post"
$arrTend[0,60] = "Shot 3pt Left Center"
$arrTend[0,61] = "Shot 3pt Right"
$arrTend[0,62] = "Shot Close Left"
$arrTend[0,63] = "Shot Close Middle"
$arrTend[0,64] = "Shot Close Right"
$arrTend[0,65] = "Shot Mid Left Center"
$arrTend[0,66] = "Shot Mid Right"
$arrTend[0,67] = "Shot Mid Right Center"
$arrTend[0,68] = "Shot Three"
$arrTend[0,69] = "Shot Three Center"
$arrTend[0,70] = "Shot Under Basket"
$arrTend[0,71] = "Spin Jumper Tendency"
$arrTend[0,72] = "Spin Layup"
$arrTend[0,73] = "Spot Up Drive"
$arrTend[0,74] = "Spot Up Shot Mid"
$arrTend[0,75] = "Standing Dunk Tendency"
$arrTend[0,76] = "Step Through Shot"
$arrTend[0,77] = "Stepback Jumper 3pt"
$arrTend[0,78] = "Stepback Jumper Mid"
$arrTend[0,79] = "Take Charge"
$arrTend[0,80] = "Touches"
$arrTend[0,81] = "Transition Pull Up 3pt"
$arrTend[0,82] = "Transition Spot Up"
$arrTend[0,83] = "Triple Threat Idle"
$arrTend[0,84] = "Triple Threat Jab Step"
$arrTend[0,85] = "Triple Threat Pump Fake"
$arrTend[0,86] = "Triple Threat Shoot"
$arrTend[0,87] = "Use Glass"
$wsTend.Range("A1").Resize(1,88).Value = $arrTend

$wsAcc = $wb.Worksheets.Item("Accessories")
$wsAcc.Range("A1:BE1").ClearContents()
$arrAcc = New-Object 'object[,]' 1,53
$arrAcc[0,0] = "Arm Frequency Paired"
$arrAcc[0,1] = "Headband"
$arrAcc[0,2] = "Headband Frequency"
$arrAcc[0,3] = "Headband Logo"
$arrAcc[0,4] = "HEADBAND_CUSTOM_COLOR"
$arrAcc[0,5] = "Knee Frequency Paired"
$arrAcc[0,6] = "Left Ankle"
$arrAcc[0,7] = "Left Ankle Item Away Color"
$arrAcc[0,8] = "Left Ankle Item Home Color"
$arrAcc[0,9] = "Left Arm Frequency"
$arrAcc[0,10] = "Left Arm Item Away Color"
$arrAcc[0,11] = "Left Arm Item Home Color"
$arrAcc[0,12] = "Left Elbow"
$arrAcc[0,13] = "Left Elbow Item Away Color"
$arrAcc[0,14] = "Left Elbow Item Home Color"
$arrAcc[0,15] = "Left Fingers"
$arrAcc[0,16] = "Left Fingers Home Color"
$arrAcc[0,17] = "Left Fingers Item Away Color"
$arrAcc[0,18] = "Left Knee"
$arrAcc[0,19] = "Left Knee Away Color"
$arrAcc[0,20] = "Left Knee Item Home Color"
$arrAcc[0,21] = "Left Leg Frequency"
$arrAcc[0,22] = "Left Leg Home Color"
$arrAcc[0,23] = "Left Leg Item Away Color"
$arrAcc[0,24] = "Left Wrist"
$arrAcc[0,25] = "Left Wrist Item Away Color"
$arrAcc[0,26] = "Left Wrist Item Home Color"
$arrAcc[0,27] = "Leg Frequency Paired"
$arrAcc[0,28] = "Right Ankle"
$arrAcc[0,29] = "Right Ankle Away Color"
$arrAcc[0,30] = "Right Ankle Item Home Color"
$arrAcc[0,31] = "Right Arm"
$arrAcc[0,32] = "Right Arm Item Away Color"
$arrAcc[0,33] = "Right Arm Item Home Color"
$arrAcc[0,34] = "Right Elbow Frequency"
$arrAcc[0,35] = "Right Elbow Home Color"
$arrAcc[0,36] = "Right Elbow Item Away Color"
$arrAcc[0,37] = "Right Fingers Frequency"
$arrAcc[0,38] = "Right Fingers Item Away Color"
$arrAcc[0,39] = "Right Fingers Item Home Color"
$arrAcc[0,40] = "Right Knee Frequency"
$arrAcc[0,41] = "Right Knee Item Away Color"
$arrAcc[0,42] = "Right Knee Item Home Color"
$arrAcc[0,43] = "Right Leg Frequency"
$arrAcc[0,44] = "Right Leg Home Color"
$arrAcc[0,45] = "Right Leg Item Away Color"
$arrAcc[0,46] = "Right Wrist"
$arrAcc[0,47] = "Right Wrist Away Color"
$arrAcc[0,48] = "Right Wrist Item Home Color"
$arrAcc[0,49] = "Shorts"
$arrAcc[0,50] = "Shorts Away Color"
$arrAcc[0,51] = "Shorts Frequency"
$arrAcc[0,52] = "Shorts Home Color"
$wsAcc.Range("A1").Resize(1,53).Value = $arrAcc

$wsSig = $wb.Worksheets.Item("Signatures")
$wsSig.Range("A1:BG1").ClearContents()
$arrSig = New-Object 'object[,]' 1,56
$arrSig[0,0] = "Aggressive Breakdown Combos"
$arrSig[0,1] = "Animation Blending"
$arrSig[0,2] = "BIG_SIZE-UP"
$arrSig[0,3] = "Chew Gum"
$arrSig[0,4] = "COMBO_MOVES"
$arrSig[0,5] = "CONTESTED"
$arrSig[0,6] = "Crossover Combos"
$arrSig[0,7] = "Dribble Pull-Up"
$arrSig[0,8] = "Dribble Style"
$arrSig[0,9] = "Dunk Emotion"
$arrSig[0,10] = "Dunk Package 10"
$arrSig[0,11] = "Dunk Package 11"
$arrSig[0,12] = "Dunk Package 12"
$arrSig[0,13] = "Dunk Package 13"
$arrSig[0,14] = "Dunk Package 14"
$arrSig[0,15] = "Dunk Package 15"
$arrSig[0,16] = "Dunk Package 2"
$arrSig[0,17] = "Dunk Package 3"
$arrSig[0,18] = "Dunk Package 4"
$arrSig[0,19] = "Dunk Package 5"
$arrSig[0,20] = "Dunk Package 6"
$arrSig[0,21] = "Dunk Package 7"
$arrSig[0,22] = "Dunk Package 8"
$arrSig[0,23] = "Dunk Package 9"
$arrSig[0,24] = "Escape Moves"
$arrSig[0,25] = "Free Throw"
$arrSig[0,26] = "Go-To Dunk Package"
$arrSig[0,27] = "Go-To Shot"
$arrSig[0,28] = "Hop Jumper"
$arrSig[0,29] = "Jumpball Ritual"
$arrSig[0,30] = "Jumpshot Base"
$arrSig[0,31] = "Jumpshot Blending"
$arrSig[0,32] = "Layup Package"
$arrSig[0,33] = "Lower/Base"
$arrSig[0,34] = "Motion Style"
$arrSig[0,35] = "Moving Behind The Back"
$arrSig[0,36] = "Moving Crossover"
$arrSig[0,37] = "Moving Hesitation"
$arrSig[0,38] = "Moving Spin"
$arrSig[0,39] = "Moving Stepback"
$arrSig[0,40] = "Pass Style"
$arrSig[0,41] = "Post Fade"
$arrSig[0,42] = "Post Hook"
$arrSig[0,43] = "Post Hop Shot"
$arrSig[0,44] = "Post Spin Shot"
$arrSig[0,45] = "Pre-Game 1"
$arrSig[0,46] = "Pre-Game 2"
$arrSig[0,47] = "Regular Breakdown Combos"
$arrSig[0,48] = "Release Timing"
$arrSig[0,49] = "Signature Size-Up"
$arrSig[0,50] = "SIGNATURE_COMBOS"
$arrSig[0,51] = "SIZE-UP_ESCAPE_PACKAGES"
$arrSig[0,52] = "Spin Jumper"
$arrSig[0,53] = "Triple Threat Style"
$arrSig[0,54] = "Upper Release 1"
$arrSig[0,55] = "Upper Release 2"
$wsSig.Range("A1").Resize(1,56).Value = $arrSig

$wsStats = $wb.Worksheets.Item("Stats")
$wsStats.Range("A1:CV1").ClearContents()
$arrStats = New-Object 'object[,]' 1,53
$arrStats[0,0] = "3PT_FIELD_GOALS_ATTEMPTED#CAREER"
$arrStats[0,1] = "3PT_FIELD_GOALS_ATTEMPTED#SEASON"
$arrStats[0,2] = "3PT_FIELD_GOALS_MADE#CAREER"
$arrStats[0,3] = "3PT_FIELD_GOALS_MADE#SEASON"
$arrStats[0,4] = "ASSISTS#CAREER"
$arrStats[0,5] = "ASSISTS#SEASON"
$arrStats[0,6] = "BLOCKS#CAREER"
$arrStats[0,7] = "BLOCKS#SEASON"
$arrStats[0,8] = "CURRENT_YEAR_STATS"
$arrStats[0,9] = "DEFENSIVE_REBOUNDS#CAREER"
$arrStats[0,10] = "DEFENSIVE_REBOUNDS#SEASON"
$arrStats[0,11] = "FIELD_GOALS_ATTEMPTED#CAREER"
$arrStats[0,12] = "FIELD_GOALS_ATTEMPTED#SEASON"
$arrStats[0,13] = "FIELD_GOALS_MADE#CAREER"
$arrStats[0,14] = "FIELD_GOALS_MADE#SEASON"
$arrStats[0,15] = "FREE_THROWS_ATTEMPTED#CAREER"
$arrStats[0,16] = "FREE_THROWS_ATTEMPTED#SEASON"
$arrStats[0,17] = "FREE_THROWS_MADE#CAREER"
$arrStats[0,18] = "FREE_THROWS_MADE#SEASON"
$arrStats[0,19] = "MINUTES_PLAYED#CAREER"
$arrStats[0,20] = "MINUTES_PLAYED#SEASON"
$arrStats[0,21] = "OFFENSIVE_REBOUNDS#CAREER"
$arrStats[0,22] = "OFFENSIVE_REBOUNDS#SEASON"
$arrStats[0,23] = "POINTS#CAREER"
$arrStats[0,24] = "POINTS#SEASON"
$arrStats[0,25] = "STATS_ID#1"
$arrStats[0,26] = "STATS_ID#2"
$arrStats[0,27] = "STATS_ID#3"
$arrStats[0,28] = "STATS_ID#4"
$arrStats[0,29] = "STATS_ID#5"
$arrStats[0,30] = "STATS_ID#6"
$arrStats[0,31] = "STATS_ID#7"
$arrStats[0,32] = "STATS_ID#8"
$arrStats[0,33] = "STATS_ID#9"
$arrStats[0,34] = "STATS_ID#10"
$arrStats[0,35] = "STATS_ID#11"
$arrStats[0,36] = "STATS_ID#12"
$arrStats[0,37] = "STATS_ID#13"
$arrStats[0,38] = "STATS_ID#14"
$arrStats[0,39] = "STATS_ID#15"
$arrStats[0,40] = "STATS_ID#16"
$arrStats[0,41] = "STATS_ID#17"
$arrStats[0,42] = "STATS_ID#18"
$arrStats[0,43] = "STATS_ID#19"
$arrStats[0,44] = "STATS_ID#20"
$arrStats[0,45] = "STATS_ID#21"
$arrStats[0,46] = "STATS_ID#22"
$arrStats[0,47] = "STATS_ID#23"
$arrStats[0,48] = "STATS_ID#24"
$arrStats[0,49] = "STEALS#CAREER"
$arrStats[0,50] = "STEALS#SEASON"
$arrStats[0,51] = "TOTAL_REBOUNDS#CAREER"
$arrStats[0,52] = "TOTAL_REBOUNDS#SEASON"
$wsStats.Range("A1").Resize(1,53).Value = $arrStats

$wsContr = $wb.Worksheets.Item("Contracts")
$wsContr.Range("A1:T1").ClearContents()
$arrContr = New-Object 'object[,]' 1,18
$arrContr[0,0] = "BIRD_YEARS"
$arrContr[0,1] = "Contract Thought"
$arrContr[0,2] = "CONTRACT_OPTION"
$arrContr[0,3] = "CONTRACT_TYPE"
$arrContr[0,4] = "CONTRACT_YEAR#1"
$arrContr[0,5] = "CONTRACT_YEAR#2"
$arrContr[0,6] = "CONTRACT_YEAR#3"
$arrContr[0,7] = "CONTRACT_YEAR#4"
$arrContr[0,8] = "CONTRACT_YEAR#5"
$arrContr[0,9] = "CONTRACT_YEAR#6"
$arrContr[0,10] = "EURO-STASH_YEARS_LEFT"
$arrContr[0,11] = "EXTENSION_LENGTH"
$arrContr[0,12] = "EXTENSION_NO_TRADE"
$arrContr[0,13] = "EXTENSION_OPTION"
$arrContr[0,14] = "NO_TRADE"
$arrContr[0,15] = "ORIGINAL_CONTRACT_YEARS"
$arrContr[0,16] = "TWO-WAY_DAYS_LEFT"
$arrContr[0,17] = "YEARS_LEFT_ON_CONTRACT"
$wsContr.Range("A1").Resize(1,18).Value = $arrContr
